$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

# Update "PERIOD TO EXPIRE" (H) and "LAST UPDATE" (I) values for rows 3-7
# to reflect the new progress-as-of date (04-Nov-2025).
# The leading apostrophe on the I-column formulas forces Excel to store the
# date-looking text as literal text (matching the original inlineStr cells)
# instead of auto-converting it into a date serial number.

$ws.Range("H3").Value = 482
$ws.Range("I3").Formula = "'04-Nov-2025"

$ws.Range("H4").Value = 35
$ws.Range("I4").Formula = "'04-Nov-2025"

$ws.Range("H5").Value = -99
$ws.Range("I5").Formula = "'04-Nov-2025"

$ws.Range("H6").Value = 286
$ws.Range("I6").Formula = "'04-Nov-2025"

$ws.Range("H7").Value = 377
$ws.Range("I7").Formula = "'04-Nov-2025"
